$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Done" for RedBlackTree row, and add a new row for Fractions
$ws.Range("B6").Value = "Done"
$ws.Range("A7").Value = "Fractions"
$ws.Range("B7").Value = "On progress"

# Resize the table to include the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:B7"))

# Extend the color-scale conditional formatting over the full column range
$fcs = $ws.Range("B2").FormatConditions
$fc = $fcs.Item(2)
$fc.ModifyAppliesToRange($ws.Range("B2:B7"))

# Update selection to match the diff
$ws.Range("B7").Select()
